$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "is low demand, leads more than triple",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "is low demand, leads to more than triple", 2
)
